# Apply "Added mode codes to the dgp_processing script" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) DSD sheet: TRANSFORMATION dimension now points at a new codelist
#    CL_COM_TRANSFORMATION instead of CL_GDP_TYPE.
# ---------------------------------------------------------------------
$dsd = $wb.Worksheets.Item("DSD")
$dsd.Range("F6").Value = "CL_COM_TRANSFORMATION"
$dsd.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) INDICATOR sheet: drop the two net-tax rows (NRTAX / RLTAX) that are
#    no longer part of the model.
# ---------------------------------------------------------------------
$indicator = $wb.Worksheets.Item("INDICATOR")
$indicator.Rows.Item(7).Delete() | Out-Null
$indicator.Rows.Item(6).Delete() | Out-Null
$indicator.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) INDUSTRY_TYPE sheet: the generic "_T"/"Gross Value Added" mode row
#    is replaced with three explicit mode codes: GDP, TAX and GWA.
# ---------------------------------------------------------------------
$industry = $wb.Worksheets.Item("INDUSTRY_TYPE")
$industry.Rows.Item(3).Insert() | Out-Null
$industry.Rows.Item(4).Insert() | Out-Null
$industry.Range("A2").Value = "GDP"
$industry.Range("B2").Value = "Gross Domestic product"
$industry.Range("A3").Value = "TAX"
$industry.Range("B3").Value = "Net Tax"
$industry.Range("A4").Value = "GWA"
$industry.Range("B4").Value = "Gross Value Added"
$industry.Range("B2").Select() | Out-Null
